$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.786.20'
$ws.Range("E2").Value = '  +4.03%  '

$ws.Range("D3").Value = '3.630.97'
$ws.Range("E3").Value = '  +3.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '632.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.11%  '

$ws.Range("D7").Value = '3.630.51'
$ws.Range("E7").Value = '  +3.24%  '

$ws.Range("E9").Value = '  +2.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.148'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.31'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.18%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.441'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000230'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.57'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.17%  '

$ws.Range("D15").Value = '4.250.06'
$ws.Range("E15").Value = '  +3.29%  '

$ws.Range("D16").Value = '3.631.04'
$ws.Range("E16").Value = '  +2.99%  '

$ws.Range("D17").Value = '69.826.04'
$ws.Range("E17").Value = '  +4.22%  '

$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("E19").Value = '  +6.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +10.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '465.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.647'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.86%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.57%  '

$ws.Range("E25").Value = '  +12.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.70%  '

$ws.Range("D27").Value = '3.779.39'
$ws.Range("E27").Value = '  +3.25%  '

$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("E29").Value = '  +14.46%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.65'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.74'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.179'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.33%  '

$ws.Range("E33").Value = '  +8.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("E35").Value = '  +6.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.67'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.82%  '

$ws.Range("D37").Value = '3.632.30'
$ws.Range("E37").Value = '  +3.42%  '

$ws.Range("E38").Value = '  +5.73%  '

$ws.Range("E39").Value = '  +13.82%  '

$ws.Range("E40").Value = '  -0.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0934'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.85%  '

$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '178.62'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.14%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '31.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +17.78%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.918'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +12.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.25%  '

$ws.Range("E49").Value = '  +2.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.84'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.270'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.06%  '
